$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "328.80" or "42.487.68").
# Force text storage so values round-trip exactly as strings, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.487.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.363.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.56%  "
$ws.Range("E13").Value = "  -6.10%  "
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.721.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.372.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.562.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.60%  "
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.06%  "
$ws.Range("E24").Value = "  -8.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.37%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0900"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0358"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.36%  "
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.227"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "68.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.47%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "115.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.02%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +38.01%  "
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.600.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.47%  "
$ws.Range("E51").Value = "  -2.77%  "
